# Run Instructions.docx - formatting and grammar update
#
# 1) Split the Title paragraph: keep "Run Instructions" as the Title, add a
#    new "By: ..." byline paragraph (with the _GoBack bookmark moved onto it)
#    directly below it, and tighten the spacing on the following Heading1
#    paragraph ("Requirements").
# 2) Merge "up-to-date) versions available." into a single run (dropping the
#    gramStart/gramEnd proofErr wrapper that was around "available").
# 3) Reword "Note that while ..." -> "While ..." and re-flow/re-split the
#    rest of that paragraph's runs, adding a gramStart/gramEnd proofErr
#    wrapper around "is".

$d = $word.ActiveDocument

function New-OpcXml([string]$bodyFragment) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/part.xml" pkg:contentType="application/xml"><pkg:xmlData>' + $bodyFragment + '</pkg:xmlData></pkg:part></pkg:package>'
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Change 1: Title paragraph split + new byline paragraph + bookmark move
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)

$frag1 = '<w:p ' + $wNs + '>' +
            '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
            '<w:r><w:t>Run Instructions</w:t></w:r>' +
        '</w:p>' +
        '<w:p ' + $wNs + '>' +
            '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
            '<w:r><w:t xml:space="preserve">By: </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>Zichen</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> Jiang (jiangz26, 001320889) and Kelvin Lin (linkk4, 001401464)</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'

$r1 = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$r1.InsertXML((New-OpcXml $frag1))

# Requirements (Heading1) paragraph now picks up tighter spacing.
$reqPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Requirements*") { $reqPara = $p; break }
}
$reqPara.SpaceBefore = 0
$reqPara.LineSpacingRule = 0   # wdLineSpaceSingle
$reqPara.LineSpacing = 12

# ---------------------------------------------------------------------
# Change 2: merge "up-to-date) versions available." into one run and
# drop the gramStart/gramEnd proofErr wrapper around "available"
# ---------------------------------------------------------------------

$notePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*up-to-date) versions available*") { $notePara = $p; break }
}

$frag2 = '<w:p ' + $wNs + '>' +
            '<w:r><w:t xml:space="preserve">Note that the operating systems and web browsers listed are the newest (most </w:t></w:r>' +
            '<w:r><w:t>up-to-date) versions available.</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> You can get the latest version by upgrading an existing version of the software.</w:t></w:r>' +
        '</w:p>'

$r2 = $d.Range($notePara.Range.Start, $notePara.Range.End)
$r2.InsertXML((New-OpcXml $frag2))

# ---------------------------------------------------------------------
# Change 3: "Note that while ..." -> "While ...", re-split runs and add
# a gramStart/gramEnd proofErr wrapper around "is"
# ---------------------------------------------------------------------

$bankPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*10 digit number is accepted*") { $bankPara = $p; break }
}

$frag3 = '<w:p ' + $wNs + '>' +
            '<w:r><w:t xml:space="preserve">While </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">any 10 digit number is accepted as a valid bank account number, the only bank account number </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">one can sign in with </w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>is</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '<w:r><w:t xml:space="preserve"> 1234567890. Implementation for the other accounts was omitted as the interface would look the same, so there would be no </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">additional </w:t></w:r>' +
            '<w:r><w:t>value added to the usability or design of the application.</w:t></w:r>' +
        '</w:p>'

$r3 = $d.Range($bankPara.Range.Start, $bankPara.Range.End)
$r3.InsertXML((New-OpcXml $frag3))

Write-Host "Done applying Run Instructions edits."
